$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellData = @"
2,G,11.64940266666667
2,H,34.948208
2,I,0.4844115508208772
2,J,0.4844115508208771
2,M,2.598166333333333
2,N,7.794499
2,O,0.3466013321552429
2,P,0.3466013321552429
2,Q,30.26708581197689
2,R,272.403772307792
2,S,0.1678976888259032
2,T,0.1678976888259032
3,G,11.64940266666667
3,H,34.948208
3,I,0.4844115508208772
3,J,0.4844115508208771
3,M,4.333403333333333
3,N,13.00021
3,O,0.5780859172985858
3,P,0.5780859172985858
3,Q,50.48156034707555
3,R,454.33404312368
3,S,0.2800314957063174
3,T,0.2800314957063172
4,G,11.64940266666667
4,H,34.948208
4,I,0.4844115508208772
4,J,0.4844115508208771
4,M,0.4692043333333333
4,N,1.407613
4,O,0.06259293136852516
4,P,0.06259293136852516
4,Q,5.46595021194489
4,R,49.193551907504
4,S,0.03032073895465201
4,T,0.030320738954652
5,G,11.64940266666667
5,H,34.948208
5,I,0.4844115508208772
5,J,0.4844115508208771
5,K,2
5,L,0.6666666666666666
5,M,0.09534933333333333
5,N,0.286048
5,O,0.01271981917764605
5,P,0.01271981917764604
5,Q,1.110762777998222
5,R,9.996865001983998
5,S,0.006161627334004656
5,T,0.006161627334004654
6,I,0.4334467773856777
6,J,0.4334467773856777
6,M,2.598166333333333
6,N,7.794499
6,O,0.3466013321552429
6,P,0.3466013321552429
6,Q,27.082696074909
6,R,243.744264674181
6,S,0.1502332304602729
6,T,0.1502332304602729
7,I,0.4334467773856777
7,J,0.4334467773856777
7,M,4.333403333333333
7,N,13.00021
7,O,0.5780859172985858
7,P,0.5780859172985858
7,Q,45.17041266410999
7,R,406.5337139769899
7,S,0.2505694779051154
7,T,0.2505694779051154
8,I,0.4334467773856777
8,J,0.4334467773856777
8,M,0.4692043333333333
8,N,1.407613
8,O,0.06259293136852516
8,P,0.06259293136852516
8,Q,4.890879461282999
8,R,44.017915151547
8,S,0.02713070438881013
8,T,0.02713070438881012
9,I,0.4334467773856777
9,J,0.4334467773856777
9,K,2
9,L,0.6666666666666666
9,M,0.09534933333333333
9,N,0.286048
9,O,0.01271981917764605
9,P,0.01271981917764604
9,Q,0.9938998063679998
9,R,8.945098257311999
9,S,0.005513364631479219
9,T,0.005513364631479219
10,G,1.677245
10,H,5.031734999999999
10,I,0.06974407828492055
10,J,0.06974407828492052
10,M,2.598166333333333
10,N,7.794499
10,O,0.3466013321552429
10,P,0.3466013321552429
10,Q,4.357761491751666
10,R,39.219853425765
10,S,0.02417339044349301
10,T,0.024173390443493
11,G,1.677245
11,H,5.031734999999999
11,I,0.06974407828492055
11,J,0.06974407828492052
11,M,4.333403333333333
11,N,13.00021
11,O,0.5780859172985858
11,P,0.5780859172985858
11,Q,7.268179073816666
11,R,65.41361166434999
11,S,0.04031806947148268
11,T,0.04031806947148266
12,G,1.677245
12,H,5.031734999999999
12,I,0.06974407828492055
12,J,0.06974407828492052
12,M,0.4692043333333333
12,N,1.407613
12,O,0.06259293136852516
12,P,0.06259293136852516
12,Q,0.7869706220616666
12,R,7.082735598554999
12,S,0.004365486305449078
12,T,0.004365486305449076
13,G,1.677245
13,H,5.031734999999999
13,I,0.06974407828492055
13,J,0.06974407828492052
13,K,2
13,L,0.6666666666666666
13,M,0.09534933333333333
13,N,0.286048
13,O,0.01271981917764605
13,P,0.01271981917764604
13,Q,0.1599241925866667
13,R,1.43931773328
13,S,0.0008871320644957795
13,T,0.0008871320644957791
14,E,2
14,F,0.6666666666666666
14,G,0.2981443333333333
14,H,0.894433
14,I,0.01239759350852466
14,J,0.01239759350852466
14,M,2.598166333333333
14,N,7.794499
14,O,0.3466013321552429
14,P,0.3466013321552429
14,Q,0.7746285693407778
14,R,6.971657124067001
14,S,0.004297022425573839
14,T,0.004297022425573839
15,E,2
15,F,0.6666666666666666
15,G,0.2981443333333333
15,H,0.894433
15,I,0.01239759350852466
15,J,0.01239759350852466
15,M,4.333403333333333
15,N,13.00021
15,O,0.5780859172985858
15,P,0.5780859172985858
15,Q,1.291979647881111
15,R,11.62781683093
15,S,0.007166874215670472
15,T,0.007166874215670471
16,E,2
16,F,0.6666666666666666
16,G,0.2981443333333333
16,H,0.894433
16,I,0.01239759350852466
16,J,0.01239759350852466
16,M,0.4692043333333333
16,N,1.407613
16,O,0.06259293136852516
16,P,0.06259293136852516
16,Q,0.1398906131587778
16,R,1.259015518429
16,S,0.0007760017196139572
16,T,0.0007760017196139571
17,E,2
17,F,0.6666666666666666
17,G,0.2981443333333333
17,H,0.894433
17,I,0.01239759350852466
17,J,0.01239759350852466
17,K,2
17,L,0.6666666666666666
17,M,0.09534933333333333
17,N,0.286048
17,O,0.01271981917764605
17,P,0.01271981917764604
17,Q,0.02842786342044444
17,R,0.255850770784
17,S,0.0001576951476663921
17,T,0.0001576951476663921
"@

$lines = $cellData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $row = $parts[0]
    $col = $parts[1]
    $val = [double]$parts[2]
    $addr = "$col$row"
    $ws.Range($addr).Value = $val
}

Write-Host "Done applying $($lines.Count) cell updates"
